$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Highlight F3 ("- African Nations Championship") with the green fill
#    used elsewhere in the sheet (same fill as style index 2).
$ws.Range("F3").Interior.Color = 11073715

# 2. Add the new "- China Cup" tournament entry in H10, matching the
#    formatting used by the other cells in that row (style index 1).
$ws.Range("G10").Copy()
$ws.Range("H10").PasteSpecial(-4122)
$ws.Range("H10").Formula = '="- China Cup"'

# 3. Update the count in H12 to include the newly added entry; the total
#    in I12 (SUM(A12:H12)) recalculates automatically.
$ws.Range("H12").Value = 9
